$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F is "dSF" - update specific rows per repulled data
$ws.Range("F2").Value = -9
$ws.Range("F4").Value = -1
$ws.Range("F6").Value = -2
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = -1
